# Clean up Authors column (column E) data for rows 2-12 on Sheet1.
# The source data had an extra space inserted after each comma-separated
# author entry (a formatting fix for previously-malformed strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 5)   # Column E = Authors
    $old = $cell.Value2
    if ($old -ne $null) {
        $new = $old -replace ',( +)', ',$1 '
        $cell.Value2 = $new
    }
}
